$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# The handback transform failed for file 3e522cf7-...: update the shared
# "Status" text for row 7 everywhere it is shown (Overview + both
# language sheets all point at the same status string).
$newStatus = "Handback transform failed"
$overview.Range("B7").Value = $newStatus
$overview.Range("C7").Value = $newStatus
$zhcn.Range("C7").Value = $newStatus
$dede.Range("C7").Value = $newStatus

# Record the handback/handoff filename mismatch error detail for each
# language's row 7 ("Error Detail" column L).
$zhcn.Range("L7").Value = "Handback file name: 2eoyftr5.exd is different with handoff file name: 3e522cf7-c6a0-4d21-8ff1-0ca05b228e09.d52225079813ff697492dcc7c95bdf4e99a0b30a.zh-cn."
$dede.Range("L7").Value = "Handback file name: 2eoyftr5.exd is different with handoff file name: 3e522cf7-c6a0-4d21-8ff1-0ca05b228e09.d52225079813ff697492dcc7c95bdf4e99a0b30a.de-de."
